$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.535.74'
$ws.Range("E2").Value = '  -1.66%  '
$ws.Range("D3").Value = '1.845.64'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.74%  '
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4243'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.19%  '
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8735'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.45%  '
$ws.Range("D13").Value = '1.859.86'
$ws.Range("E13").Value = '  -5.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.331'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.507'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06908'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '79.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008860'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("E21").Value = '  -2.50%  '
$ws.Range("D22").Value = '27.542.72'
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.985'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("E24").Value = '  -5.36%  '
$ws.Range("D25").Value = '2.080.47'
$ws.Range("E25").Value = '  -4.54%  '
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.90'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '121.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.238'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.879'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08871'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7625'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.551'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.947'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("E36").Value = '  -6.67%  '
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.090'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05340'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01931'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.810'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5092'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.853'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1644'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.259'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06528'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4738'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.19%  '
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.620'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.13%  '
